# "Apr 2020 to May 2020 - Cord Cutting Wizard - Change Log.xlsx"
# The change-reason phrases in column E were reworded/re-titled. Re-apply the
# new wording to every row that used the old wording (the underlying meaning
# of each row is unchanged - only the label text was rewritten/retitled).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Apr 2020 to May 2020")

# Old label text -> new label text
$rewordings = @{
    "Old Network removed from database in May 2020"                  = "Network Removed from Database";
    "Network removed from base Service in May 2020"                  = "Network Removed from Base Service";
    "Network removed from Add-On Service in May 2020"                = "Network Removed from Add-On Package";
    "Network moved from base Service to Add-On Service in May 2020"  = "Network Moved from Base Service to Add-On Package";
    "Network added to base Service in May 2020"                      = "Network Added to Base Service";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $old = $cell.Value()
    if ($rewordings.ContainsKey($old)) {
        $cell.Value = $rewordings[$old]
    }
}

# Re-establish the selection shown when the file was saved
$ws.Activate()
$ws.Range("E2:E26").Select()
